$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: sets a cell's value forcing a Text interpretation, guarding against
# values that look numeric or date-like, which Excel's COM layer would
# otherwise silently auto-convert (losing leading zeros / becoming date serials),
# while leaving the original OOXML cell type as text (inline string), matching
# the source data model used throughout this worksheet.
function Set-TextCell($cell, [string]$val) {
    if ($val -match '^[+-]?[0-9]+(\.[0-9]+)?$' -or $val -match '^[0-9]{1,2}/[0-9]{1,2}/[0-9]{4}$') {
        $cell.Value2 = "'" + $val
    } else {
        $cell.Value2 = $val
    }
}

# 1) Delete row 2 (Caso 4238, GUATEMALA 5527) - all subsequent rows shift up by one
$ws.Rows.Item(2).Delete()

# 2) Append 10 new case rows at the end (rows 31-40)

# Row 31: Caso 7640 - CIUDAD DE LA PAZ 258
Set-TextCell ($ws.Cells.Item(31, 1)) '7640'
Set-TextCell ($ws.Cells.Item(31, 2)) '10/28/2025'
Set-TextCell ($ws.Cells.Item(31, 3)) 'CIUDAD DE LA PAZ 258'
$ws.Cells.Item(31, 4).Value2 = 14
Set-TextCell ($ws.Cells.Item(31, 5)) '01183611'
Set-TextCell ($ws.Cells.Item(31, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(31, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(31, 8)) 'Cable en panza o cortados'
$ws.Cells.Item(31, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(31, 10)) '{"direccionesNormalizadas": [{"altura": 258, "cod_calle": 3128, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.440983", "y": "-34.574792"}, "direccion": "CIUDAD DE LA PAZ 258, CABA", "nombre_calle": "CIUDAD DE LA PAZ", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(31, 11).Value2 = -58.440983
$ws.Cells.Item(31, 12).Value2 = -34.574792
Set-TextCell ($ws.Cells.Item(31, 13)) 'Palermo'
Set-TextCell ($ws.Cells.Item(31, 14)) 'Capital Sur'

# Row 32: Caso 7680 - LARREA 811
Set-TextCell ($ws.Cells.Item(32, 1)) '7680'
Set-TextCell ($ws.Cells.Item(32, 2)) '10/28/2025'
Set-TextCell ($ws.Cells.Item(32, 3)) 'LARREA 811'
$ws.Cells.Item(32, 4).Value2 = 2
Set-TextCell ($ws.Cells.Item(32, 5)) '01196663'
Set-TextCell ($ws.Cells.Item(32, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(32, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(32, 8)) 'Cable en panza'
$ws.Cells.Item(32, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(32, 10)) '{"direccionesNormalizadas": [{"altura": 811, "cod_calle": 12065, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.402416", "y": "-34.598847"}, "direccion": "LARREA 811, CABA", "nombre_calle": "LARREA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(32, 11).Value2 = -58.402416
$ws.Cells.Item(32, 12).Value2 = -34.598847
Set-TextCell ($ws.Cells.Item(32, 13)) 'Recoleta'
Set-TextCell ($ws.Cells.Item(32, 14)) 'Capital Sur'

# Row 33: Caso 3447 - AZOPARDO 1071
Set-TextCell ($ws.Cells.Item(33, 1)) '3447'
Set-TextCell ($ws.Cells.Item(33, 2)) '10/28/2025'
Set-TextCell ($ws.Cells.Item(33, 3)) 'AZOPARDO 1071'
$ws.Cells.Item(33, 4).Value2 = 1
Set-TextCell ($ws.Cells.Item(33, 5)) '01196692'
Set-TextCell ($ws.Cells.Item(33, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(33, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(33, 8)) 'Tendido a baja altura'
$ws.Cells.Item(33, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(33, 10)) '{"direccionesNormalizadas": [{"altura": 1071, "cod_calle": 1149, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.367154", "y": "-34.620075"}, "direccion": "AZOPARDO 1071, CABA", "nombre_calle": "AZOPARDO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(33, 11).Value2 = -58.367154
$ws.Cells.Item(33, 12).Value2 = -34.620075
Set-TextCell ($ws.Cells.Item(33, 13)) 'San Telmo'
Set-TextCell ($ws.Cells.Item(33, 14)) 'Capital Sur'

# Row 34: Caso 3487 - ACOSTA, MARIANO AV. 171
Set-TextCell ($ws.Cells.Item(34, 1)) '3487'
Set-TextCell ($ws.Cells.Item(34, 2)) '10/28/2025'
Set-TextCell ($ws.Cells.Item(34, 3)) 'ACOSTA, MARIANO AV. 171'
$ws.Cells.Item(34, 4).Value2 = 10
Set-TextCell ($ws.Cells.Item(34, 5)) '01188247'
Set-TextCell ($ws.Cells.Item(34, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(34, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(34, 8)) 'Cable en panza'
$ws.Cells.Item(34, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(34, 10)) '{"direccionesNormalizadas": [{"altura": 171, "cod_calle": 1006, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.480436", "y": "-34.635569"}, "direccion": "ACOSTA, MARIANO AV. 171, CABA", "nombre_calle": "ACOSTA, MARIANO AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(34, 11).Value2 = -58.480436
$ws.Cells.Item(34, 12).Value2 = -34.635569
Set-TextCell ($ws.Cells.Item(34, 13)) 'Devoto'
Set-TextCell ($ws.Cells.Item(34, 14)) 'Capital Norte'

# Row 35: Caso 7696 - CORDOBA AV. 1776
Set-TextCell ($ws.Cells.Item(35, 1)) '7696'
Set-TextCell ($ws.Cells.Item(35, 2)) '10/29/2025'
Set-TextCell ($ws.Cells.Item(35, 3)) 'CORDOBA AV. 1776'
$ws.Cells.Item(35, 4).Value2 = 1
Set-TextCell ($ws.Cells.Item(35, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(35, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(35, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(35, 8)) 'Cable en panza'
$ws.Cells.Item(35, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(35, 10)) '{"direccionesNormalizadas": [{"altura": 1776, "cod_calle": 3165, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.392498", "y": "-34.599695"}, "direccion": "CORDOBA AV. 1776, CABA", "nombre_calle": "CORDOBA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(35, 11).Value2 = -58.392498
$ws.Cells.Item(35, 12).Value2 = -34.599695
Set-TextCell ($ws.Cells.Item(35, 13)) 'San Telmo'
Set-TextCell ($ws.Cells.Item(35, 14)) 'Capital Sur'

# Row 36: Caso 7698 - 24 DE NOVIEMBRE 151
Set-TextCell ($ws.Cells.Item(36, 1)) '7698'
Set-TextCell ($ws.Cells.Item(36, 2)) '10/30/2025'
Set-TextCell ($ws.Cells.Item(36, 3)) '24 DE NOVIEMBRE 151'
$ws.Cells.Item(36, 4).Value2 = 3
Set-TextCell ($ws.Cells.Item(36, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(36, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(36, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(36, 8)) 'Cable en panza'
$ws.Cells.Item(36, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(36, 10)) '{"direccionesNormalizadas": [{"altura": 151, "cod_calle": 23025, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.411626", "y": "-34.612418"}, "direccion": "24 DE NOVIEMBRE 151, CABA", "nombre_calle": "24 DE NOVIEMBRE", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(36, 11).Value2 = -58.411626
$ws.Cells.Item(36, 12).Value2 = -34.612418
Set-TextCell ($ws.Cells.Item(36, 13)) 'Almagro'
Set-TextCell ($ws.Cells.Item(36, 14)) 'Capital Sur'

# Row 37: Caso 7711 - LARREA 608
Set-TextCell ($ws.Cells.Item(37, 1)) '7711'
Set-TextCell ($ws.Cells.Item(37, 2)) '10/30/2025'
Set-TextCell ($ws.Cells.Item(37, 3)) 'LARREA 608'
$ws.Cells.Item(37, 4).Value2 = 3
Set-TextCell ($ws.Cells.Item(37, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(37, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(37, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(37, 8)) 'Cable en panza cables cortados'
$ws.Cells.Item(37, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(37, 10)) '{"direccionesNormalizadas": [{"altura": 608, "cod_calle": 12065, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.402363", "y": "-34.601960"}, "direccion": "LARREA 608, CABA", "nombre_calle": "LARREA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(37, 11).Value2 = -58.402363
$ws.Cells.Item(37, 12).Value2 = -34.60196
Set-TextCell ($ws.Cells.Item(37, 13)) 'Almagro'
Set-TextCell ($ws.Cells.Item(37, 14)) 'Capital Sur'

# Row 38: Caso 7722 - CONCORDIA 1401
Set-TextCell ($ws.Cells.Item(38, 1)) '7722'
Set-TextCell ($ws.Cells.Item(38, 2)) '10/30/2025'
Set-TextCell ($ws.Cells.Item(38, 3)) 'CONCORDIA 1401'
$ws.Cells.Item(38, 4).Value2 = 11
Set-TextCell ($ws.Cells.Item(38, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(38, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(38, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(38, 8)) 'Cable en panza'
$ws.Cells.Item(38, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(38, 10)) '{"direccionesNormalizadas": [{"altura": 1401, "cod_calle": 3151, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.484849", "y": "-34.618754"}, "direccion": "CONCORDIA 1401, CABA", "nombre_calle": "CONCORDIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(38, 11).Value2 = -58.484849
$ws.Cells.Item(38, 12).Value2 = -34.618754
Set-TextCell ($ws.Cells.Item(38, 13)) 'Devoto'
Set-TextCell ($ws.Cells.Item(38, 14)) 'Capital Norte'

# Row 39: Caso 7729 - SAN JOSE 343
Set-TextCell ($ws.Cells.Item(39, 1)) '7729'
Set-TextCell ($ws.Cells.Item(39, 2)) '10/30/2025'
Set-TextCell ($ws.Cells.Item(39, 3)) 'SAN JOSE 343'
$ws.Cells.Item(39, 4).Value2 = 1
Set-TextCell ($ws.Cells.Item(39, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(39, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(39, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(39, 8)) 'Cable en panza'
$ws.Cells.Item(39, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(39, 10)) '{"direccionesNormalizadas": [{"altura": 343, "cod_calle": 20038, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.386186", "y": "-34.612745"}, "direccion": "SAN JOSE 343, CABA", "nombre_calle": "SAN JOSE", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}, {"altura": 343, "cod_calle": 20039, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.437866", "y": "-34.622981"}, "direccion": "SAN JOSE DE CALASANZ 343, CABA", "nombre_calle": "SAN JOSE DE CALASANZ", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(39, 11).Value2 = -58.386186
$ws.Cells.Item(39, 12).Value2 = -34.612745
Set-TextCell ($ws.Cells.Item(39, 13)) 'San Telmo'
Set-TextCell ($ws.Cells.Item(39, 14)) 'Capital Sur'

# Row 40: Caso 7740 - DEL LIBERTADOR AV. 8620
Set-TextCell ($ws.Cells.Item(40, 1)) '7740'
Set-TextCell ($ws.Cells.Item(40, 2)) '10/30/2025'
Set-TextCell ($ws.Cells.Item(40, 3)) 'DEL LIBERTADOR AV. 8620'
$ws.Cells.Item(40, 4).Value2 = 13
Set-TextCell ($ws.Cells.Item(40, 5)) 'Pendiente ADM'
Set-TextCell ($ws.Cells.Item(40, 6)) 'Optical Power'
Set-TextCell ($ws.Cells.Item(40, 7)) 'Pendiente'
Set-TextCell ($ws.Cells.Item(40, 8)) 'Tendido a baja altura'
$ws.Cells.Item(40, 9).Value2 = 1
Set-TextCell ($ws.Cells.Item(40, 10)) '{"direccionesNormalizadas": [{"altura": 8620, "cod_calle": 12107, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.466588", "y": "-34.536500"}, "direccion": "DEL LIBERTADOR AV. 8620, CABA", "nombre_calle": "DEL LIBERTADOR AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(40, 11).Value2 = -58.466588
$ws.Cells.Item(40, 12).Value2 = -34.5365
Set-TextCell ($ws.Cells.Item(40, 13)) 'Saavedra'
Set-TextCell ($ws.Cells.Item(40, 14)) 'Capital Norte'
